{"js": "// Update the date line and the 25 division problems in the table.\nconst replacements = [\n  [\"2025-05-01 Thursday\", \"2025-05-02 Friday\"],\n  [\"77\u00f75=\", \"45\u00f76=\"],\n  [\"35\u00f78=\", \"41\u00f73=\"],\n  [\"59\u00f72=\", \"49\u00f75=\"],\n  [\"40\u00f79=\", \"86\u00f79=\"],\n  [\"96\u00f78=\", \"68\u00f77=\"],\n  [\"77\u00f74=\", \"80\u00f74=\"],\n  [\"40\u00f72=\", \"71\u00f78=\"],\n  [\"80\u00f73=\", \"43\u00f74=\"],\n  [\"71\u00f72=\", \"78\u00f78=\"],\n  [\"25\u00f75=\", \"43\u00f78=\"],\n  [\"47\u00f79=\", \"75\u00f75=\"],\n  [\"63\u00f72=\", \"69\u00f79=\"],\n  [\"58\u00f79=\", \"85\u00f73=\"],\n  [\"41\u00f79=\", \"14\u00f79=\"],\n  [\"63\u00f76=\", \"61\u00f73=\"],\n  [\"31\u00f73=\", \"84\u00f78=\"],\n  [\"52\u00f79=\", \"88\u00f76=\"],\n  [\"37\u00f74=\", \"62\u00f76=\"],\n  [\"61\u00f79=\", \"93\u00f74=\"],\n  [\"88\u00f78=\", \"76\u00f78=\"],\n  [\"30\u00f76=\", \"86\u00f73=\"],\n  [\"16\u00f79=\", \"54\u00f75=\"],\n  [\"33\u00f74=\", \"48\u00f76=\"],\n  [\"90\u00f74=\", \"55\u00f74=\"],\n  [\"46\u00f74=\", \"89\u00f72=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date line and the 25 division problems in the table.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-05-01 Thursday\", \"2025-05-02 Friday\"),\n    @(\"77\u00f75=\", \"45\u00f76=\"),\n    @(\"35\u00f78=\", \"41\u00f73=\"),\n    @(\"59\u00f72=\", \"49\u00f75=\"),\n    @(\"40\u00f79=\", \"86\u00f79=\"),\n    @(\"96\u00f78=\", \"68\u00f77=\"),\n    @(\"77\u00f74=\", \"80\u00f74=\"),\n    @(\"40\u00f72=\", \"71\u00f78=\"),\n    @(\"80\u00f73=\", \"43\u00f74=\"),\n    @(\"71\u00f72=\", \"78\u00f78=\"),\n    @(\"25\u00f75=\", \"43\u00f78=\"),\n    @(\"47\u00f79=\", \"75\u00f75=\"),\n    @(\"63\u00f72=\", \"69\u00f79=\"),\n    @(\"58\u00f79=\", \"85\u00f73=\"),\n    @(\"41\u00f79=\", \"14\u00f79=\"),\n    @(\"63\u00f76=\", \"61\u00f73=\"),\n    @(\"31\u00f73=\", \"84\u00f78=\"),\n    @(\"52\u00f79=\", \"88\u00f76=\"),\n    @(\"37\u00f74=\", \"62\u00f76=\"),\n    @(\"61\u00f79=\", \"93\u00f74=\"),\n    @(\"88\u00f78=\", \"76\u00f78=\"),\n    @(\"30\u00f76=\", \"86\u00f73=\"),\n    @(\"16\u00f79=\", \"54\u00f75=\"),\n    @(\"33\u00f74=\", \"48\u00f76=\"),\n    @(\"90\u00f74=\", \"55\u00f74=\"),\n    @(\"46\u00f74=\", \"89\u00f72=\")\n)\n\nforeach ($pair in $replacements) {\n    $find = $pair[0]\n    $replace = $pair[1]\n\n    $range = $d.Content\n    $range.Find.ClearFormatting()\n    $range.Find.Replacement.ClearFormatting()\n    $range.Find.Text = $find\n    $range.Find.Replacement.Text = $replace\n    $range.Find.Forward = $true\n    $range.Find.Wrap = 1  # wdFindContinue\n    $range.Find.MatchCase = $true\n    $range.Find.MatchWholeWord = $false\n    $range.Find.MatchWildcards = $false\n    $range.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2) | Out-Null\n}\n"}
